$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (shifts existing rows 80..171 down to 81..172)
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new data record
$ws.Range("A80").Value = 3
$ws.Range("B80").Value = "Femacal de La Calera"
$ws.Range("C80").Value = "Coquimbo"
$ws.Range("D80").Value = 44664
$ws.Range("E80").Value = 5
$ws.Range("F80").Value = 100112052
$ws.Range("G80").Value = "Albahaca"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 110
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = 3727
$ws.Range("N80").Value = "`$/docena de matas"
$ws.Range("O80").Value = "Provincia de Quillota"
$ws.Range("P80").Value = 621
$ws.Range("Q80").Value = 6
$ws.Range("R80").Value = "Hortaliza"
